# "Generate Report for Handback"
#
# 26b6cf91-7de0-4f04-9e05-f2f12ab2b611 has now also been handed back
# (in sync with en-US), same as c48d1ac2-dd35-449e-b7df-405ee16c0b4c
# already was. The two files swap display order (26b6cf91 now listed
# first / row 2, c48d1ac2 second / row 3) on every sheet, and the
# 26b6cf91 rows gain their own real "Latest Handback File" /
# "Latest Handback DateTime" values instead of the placeholder
# "Ready for handoff" status.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-39-12 16:39:19"

$ws.Range("A3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-37-12 16:37:53"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.1ba2234374188bd1f9b5e1320a5980e88f9fcd21.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-12 16:39:16"
$ws.Range("F2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.md"
$ws.Range("G2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.1ba2234374188bd1f9b5e1320a5980e88f9fcd21.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-12 16:39:32"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.5582d30980d297f5b7910f0ff5e1d3cb573b37a6.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-12 16:37:12"
$ws.Range("F3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.md"
$ws.Range("G3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.5582d30980d297f5b7910f0ff5e1d3cb573b37a6.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-12 16:38:49"
$ws.Range("I3").Value = "Include"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.1ba2234374188bd1f9b5e1320a5980e88f9fcd21.de-de.xlf"
$ws.Range("E2").Value = "2016-03-12 16:39:19"
$ws.Range("F2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.md"
$ws.Range("G2").Value = "26b6cf91-7de0-4f04-9e05-f2f12ab2b611.1ba2234374188bd1f9b5e1320a5980e88f9fcd21.de-de.xlf"
$ws.Range("H2").Value = "2016-03-12 16:39:37"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.5582d30980d297f5b7910f0ff5e1d3cb573b37a6.de-de.xlf"
$ws.Range("E3").Value = "2016-03-12 16:37:53"
$ws.Range("F3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.md"
$ws.Range("G3").Value = "c48d1ac2-dd35-449e-b7df-405ee16c0b4c.5582d30980d297f5b7910f0ff5e1d3cb573b37a6.de-de.xlf"
$ws.Range("H3").Value = "2016-03-12 16:38:54"
$ws.Range("I3").Value = "Include"
